# Data-cleanup pass ahead of the SQL export: the "IPs" column (C) has a
# handful of rows where the inflation-pressure code was mistakenly entered
# with a leading "P" ("P2" instead of "2"), and the "Test Velocity" column
# (H) has a handful of "Rolling Resistance" rows where "10" should read "9".
# Both columns store their data as text, so we re-enter the corrected
# values with a leading apostrophe (forces text, matching every other
# value already in these columns) and then reset the cell style back to
# "Normal" so no stray number-format/quote-prefix styling is left behind.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("IPs"): "P2" -> "2"
$cRanges = @("C9:C26", "C39:C56", "C63:C86", "C102:C116", "C132:C134", "C138:C140", "C147:C149", "C155:C156", "C171:C182")
foreach ($addr in $cRanges) {
    $rng = $ws.Range($addr)
    $rng.Value = "'2"
    $rng.Style = "Normal"
}

# Column H ("Test Velocity"): "10" -> "9" (Rolling Resistance rows)
$hRanges = @("H168:H170", "H180:H182", "H192:H194")
foreach ($addr in $hRanges) {
    $rng = $ws.Range($addr)
    $rng.Value = "'9"
    $rng.Style = "Normal"
}
